$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 963-1022 (appended block of stock-screener data)
$rowsData = @(
    @{ Row=963; "A"="Buying Opportunity"; "B"="support Zone"; "C"="long buildup"; "D"="Short buildup"; "E"="FII ENTERING" };
    @{ Row=964; "A"="AIRAN"; "B"="HEIDELBERG"; "C"="ABFRL"; "E"="COROMANDEL"; "F"=29.5; "G"=217.59; "H"=330.5; "J"=1530.15 };
    @{ Row=965; "A"="ALPHAETF"; "B"="MUFIN"; "E"="MOTHERSON"; "F"=27.89; "G"=125.21; "J"=180.32 };
    @{ Row=966; "A"="ANANDRATHI"; "B"="PANSARI"; "F"=4014.9; "G"=97.39 };
    @{ Row=967; "A"="ANIKINDS"; "F"=53.54 };
    @{ Row=968; "A"="ASHOKA"; "F"=228.07 };
    @{ Row=969; "A"="ATL"; "F"=54.91 };
    @{ Row=970; "A"="AXISCADES"; "F"=617.55 };
    @{ Row=971; "A"="BAJAJELEC"; "F"=1066.65 };
    @{ Row=972; "A"="BALMLAWRIE"; "F"=271.7 };
    @{ Row=973; "A"="BFSI"; "F"=23.04 };
    @{ Row=974; "A"="CANFINHOME"; "F"=889.25 };
    @{ Row=975; "A"="CARBORUNIV"; "F"=1783.05 };
    @{ Row=976; "A"="CERA"; "F"=8078.55 };
    @{ Row=977; "A"="COCHINSHIP"; "F"=2320.8 };
    @{ Row=978; "A"="CONSUMBEES"; "F"=124.01 };
    @{ Row=979; "A"="COROMANDEL"; "F"=1530.15 };
    @{ Row=980; "A"="DATAPATTNS"; "F"=3088.1 };
    @{ Row=981; "A"="ESG"; "F"=39.68 };
    @{ Row=982; "A"="FACT"; "F"=866.8 };
    @{ Row=983; "A"="GIPCL"; "F"=252.53 };
    @{ Row=984; "A"="GMRINFRA"; "F"=97.84 };
    @{ Row=985; "A"="GPIL"; "F"=1088.9 };
    @{ Row=986; "A"="GPPL"; "F"=202.21 };
    @{ Row=987; "A"="GRSE"; "F"=1847.55 };
    @{ Row=988; "A"="GULFPETRO"; "F"=68.71 };
    @{ Row=989; "A"="HAL"; "F"=5533.45 };
    @{ Row=990; "A"="HDFCBANK"; "F"=1607.8 };
    @{ Row=991; "A"="HDFCMOMENT"; "F"=36.69 };
    @{ Row=992; "A"="HDFCNEXT50"; "F"=72.96 };
    @{ Row=993; "A"="HERCULES"; "F"=556.05 };
    @{ Row=994; "A"="HFCL"; "F"=124.19 };
    @{ Row=995; "A"="IDEAFORGE"; "F"=830.7 };
    @{ Row=996; "A"="IGPL"; "F"=573.25 };
    @{ Row=997; "A"="IMPAL"; "F"=1114.6 };
    @{ Row=998; "A"="JSWENERGY"; "F"=706.7 };
    @{ Row=999; "A"="KAPSTON"; "F"=376.5 };
    @{ Row=1000; "A"="KAYNES"; "F"=3894.45 };
    @{ Row=1001; "A"="KBCGLOBAL"; "F"=1.9 };
    @{ Row=1002; "A"="KICL"; "F"=4607.9 };
    @{ Row=1003; "A"="M&M"; "F"=2961.9 };
    @{ Row=1004; "A"="MIDHANI"; "F"=461.4 };
    @{ Row=1005; "A"="MOMENTUM"; "F"=36.69 };
    @{ Row=1006; "A"="MOMOMENTUM"; "F"=73.4 };
    @{ Row=1007; "A"="MONIFTY500"; "F"=22.57 };
    @{ Row=1008; "A"="MOTHERSON"; "F"=180.32 };
    @{ Row=1009; "A"="MTARTECH"; "F"=1911.45 };
    @{ Row=1010; "A"="NACLIND"; "F"=74.71 };
    @{ Row=1011; "A"="NDLVENTURE"; "F"=101.45 };
    @{ Row=1012; "A"="NETF"; "F"=250.35 };
    @{ Row=1013; "A"="NFL"; "F"=120.44 };
    @{ Row=1014; "A"="NIF100BEES"; "F"=257.94 };
    @{ Row=1015; "A"="NIFTYQLITY"; "F"=21.12 };
    @{ Row=1016; "A"="PHOENIXLTD"; "F"=3804.55 };
    @{ Row=1017; "A"="PILANIINVS"; "F"=3944.1 };
    @{ Row=1018; "A"="PNBHOUSING"; "F"=864.8 };
    @{ Row=1019; "A"="PRAKASH"; "F"=189.38 };
    @{ Row=1020; "A"="QUICKHEAL"; "F"=524.9 };
    @{ Row=1021; "A"="RADIANTCMS"; "F"=83.23 };
    @{ Row=1022; "A"="19/06/2024" }
)

foreach ($rowData in $rowsData) {
    $r = $rowData.Row
    foreach ($col in $rowData.Keys) {
        if ($col -ne "Row") {
            $ws.Range("$col$r").Value = $rowData[$col]
        }
    }
}
